# Applies the commit: adds two new weekly price records (rows) for
# Frambuesa / Mercado Mayorista Lo Valledor de Santiago, inserted right
# before the existing row that used to be row 112 (date 2021-02-09 /
# serial 44260). This pushes all subsequent rows down by two, growing
# the sheet from A1:T223 to A1:T225.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the old row 112 position, shifting the rest
# of the table (old rows 112..223) down to 114..225.
$ws.Rows("112:113").Insert()

# --- New row 112 --------------------------------------------------------
$ws.Range("A112").Value = 6
$ws.Range("B112").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C112").Value = "Metropolitana"
$ws.Range("D112").Value = 44902
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100101
$ws.Range("H112").Value = "Berries"
$ws.Range("I112").Value = 100101004
$ws.Range("J112").Value = "Frambuesa"
$ws.Range("K112").Value = "Sin especificar"
$ws.Range("L112").Value = "Especial"
$ws.Range("M112").Value = 250
$ws.Range("N112").Value = 8000
$ws.Range("O112").Value = 8000
$ws.Range("P112").Value = 8000
$ws.Range("Q112").Value = "`$/bandeja 2 kilos"
$ws.Range("R112").Value = "Región de O'Higgins"
$ws.Range("S112").Value = 4000
$ws.Range("T112").Value = 2

# --- New row 113 --------------------------------------------------------
$ws.Range("A113").Value = 6
$ws.Range("B113").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C113").Value = "Metropolitana"
$ws.Range("D113").Value = 44902
$ws.Range("E113").Value = 13
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100101
$ws.Range("H113").Value = "Berries"
$ws.Range("I113").Value = 100101004
$ws.Range("J113").Value = "Frambuesa"
$ws.Range("K113").Value = "Sin especificar"
$ws.Range("L113").Value = "Primera"
$ws.Range("M113").Value = 350
$ws.Range("N113").Value = 7000
$ws.Range("O113").Value = 7000
$ws.Range("P113").Value = 7000
$ws.Range("Q113").Value = "`$/bandeja 2 kilos"
$ws.Range("R113").Value = "Región del Maule"
$ws.Range("S113").Value = 3500
$ws.Range("T113").Value = 2
